$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.027592182159424
$ws.Range("B1").Value = 1.367025971412659
$ws.Range("D1").Value = 1.707218766212463
$ws.Range("E1").Value = 1.025221943855286
